$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 19.68388059701489
$ws.Range("B2").Value = 0.7400000000000003
$ws.Range("C2").Value = 52.80666666666666
$ws.Range("A3").Value = 4.83369817578772
$ws.Range("B3").Value = 0.09333333333333335
$ws.Range("C3").Value = 15.27333333333334
$ws.Range("A4").Value = 6.794195688225533
$ws.Range("B4").Value = 0.1333333333333334
$ws.Range("C4").Value = 19.42
$ws.Range("A5").Value = 6.200862354892194
$ws.Range("B5").Value = 0.07333333333333332
$ws.Range("C5").Value = 18.64666666666666
$ws.Range("A6").Value = 12.32371475953564
$ws.Range("B6").Value = 0.5533333333333333
$ws.Range("C6").Value = 36.36666666666668
$ws.Range("A7").Value = 18.49568822553895
$ws.Range("B7").Value = 0.4866666666666667
$ws.Range("C7").Value = 52.12666666666668
$ws.Range("A8").Value = 25.75442786069642
$ws.Range("B8").Value = 0.8666666666666661
$ws.Range("C8").Value = 66.59333333333336
$ws.Range("A9").Value = 11.32527363184078
$ws.Range("B9").Value = 0.4333333333333335
$ws.Range("C9").Value = 31.92000000000001
$ws.Range("A10").Value = 24.39917081260354
$ws.Range("B10").Value = 0.5
$ws.Range("C10").Value = 67.43333333333335
$ws.Range("A11").Value = 22.48175787728016
$ws.Range("B11").Value = 0.9466666666666667
$ws.Range("C11").Value = 61.21333333333333
$ws.Range("A12").Value = 18.34796019900498
$ws.Range("B12").Value = 0.4733333333333336
$ws.Range("C12").Value = 46.84666666666665
$ws.Range("A13").Value = 25.098839137645
$ws.Range("B13").Value = 1.08
$ws.Range("C13").Value = 68.24666666666667
$ws.Range("A14").Value = 24.02736318407949
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 63.94
$ws.Range("A15").Value = 7.255058043117734
$ws.Range("B15").Value = 0.1133333333333334
$ws.Range("C15").Value = 22.96
$ws.Range("A16").Value = 12.28839137645106
$ws.Range("B16").Value = 0.4066666666666666
$ws.Range("C16").Value = 35.18000000000001
$ws.Range("A17").Value = 16.64600331674958
$ws.Range("B17").Value = 0.7266666666666668
$ws.Range("C17").Value = 46.64666666666668
$ws.Range("A18").Value = 6.619502487562184
$ws.Range("B18").Value = 0.1266666666666667
$ws.Range("C18").Value = 18.97333333333333
$ws.Range("A19").Value = 23.15160862354889
$ws.Range("B19").Value = 0.6333333333333333
$ws.Range("C19").Value = 63.81999999999998
$ws.Range("A20").Value = 14.1372470978441
$ws.Range("B20").Value = 0.4266666666666665
$ws.Range("C20").Value = 38.47333333333334
$ws.Range("A21").Value = 9.384112769485901
$ws.Range("B21").Value = 0.1066666666666667
$ws.Range("C21").Value = 27.84000000000001
$ws.Range("A22").Value = 24.58494195688216
$ws.Range("B22").Value = 0.9800000000000001
$ws.Range("C22").Value = 67.49999999999999
$ws.Range("A23").Value = 8.063250414593691
$ws.Range("B23").Value = 0.1666666666666667
$ws.Range("C23").Value = 23.74
$ws.Range("A24").Value = 13.40643449419568
$ws.Range("B24").Value = 0.3333333333333333
$ws.Range("C24").Value = 37.20000000000001
$ws.Range("A25").Value = 10.43993366500828
$ws.Range("B25").Value = 0.2200000000000001
$ws.Range("C25").Value = 29.35333333333333
$ws.Range("A26").Value = 8.951641791044766
$ws.Range("B26").Value = 0.3666666666666668
$ws.Range("C26").Value = 29.08000000000001
$ws.Range("A27").Value = 25.28693200663339
$ws.Range("B27").Value = 0.9600000000000005
$ws.Range("C27").Value = 69.79333333333334
$ws.Range("A28").Value = 20.19021558872293
$ws.Range("B28").Value = 0.6533333333333334
$ws.Range("C28").Value = 52.04000000000001
$ws.Range("A29").Value = 8.251575456053054
$ws.Range("B29").Value = 0.1466666666666667
$ws.Range("C29").Value = 24.1
$ws.Range("A30").Value = 15.61124378109452
$ws.Range("B30").Value = 0.64
$ws.Range("C30").Value = 44.21333333333334
$ws.Range("A31").Value = 24.21442786069639
$ws.Range("B31").Value = 0.8533333333333331
$ws.Range("C31").Value = 66.81333333333332
$ws.Range("A32").Value = 24.89940298507456
$ws.Range("B32").Value = 0.9400000000000001
$ws.Range("C32").Value = 66.80000000000001
$ws.Range("A33").Value = 24.75164179104467
$ws.Range("B33").Value = 0.9266666666666669
$ws.Range("C33").Value = 64.81333333333333
$ws.Range("A34").Value = 21.72245439469313
$ws.Range("B34").Value = 0.78
$ws.Range("C34").Value = 59.11333333333335
$ws.Range("A35").Value = 13.72905472636815
$ws.Range("B35").Value = 0.2933333333333333
$ws.Range("C35").Value = 37.98666666666667
$ws.Range("A36").Value = 19.10686567164178
$ws.Range("B36").Value = 0.6733333333333333
$ws.Range("C36").Value = 50.49333333333334
$ws.Range("A37").Value = 10.73986733001658
$ws.Range("B37").Value = 0.2333333333333333
$ws.Range("C37").Value = 30.64000000000001
$ws.Range("A38").Value = 24.19927031509116
$ws.Range("B38").Value = 1.12
$ws.Range("C38").Value = 66.8133333333333
$ws.Range("A39").Value = 18.25442786069652
$ws.Range("B39").Value = 0.5599999999999998
$ws.Range("C39").Value = 47.74666666666666
$ws.Range("A40").Value = 8.33296849087893
$ws.Range("B40").Value = 0.1800000000000001
$ws.Range("C40").Value = 25.88666666666667
$ws.Range("A41").Value = 17.82922056384741
$ws.Range("B41").Value = 0.5800000000000001
$ws.Range("C41").Value = 45.61333333333334
$ws.Range("A42").Value = 22.45432835820887
$ws.Range("B42").Value = 0.9599999999999997
$ws.Range("C42").Value = 62.42666666666666
$ws.Range("A43").Value = 15.80056384742951
$ws.Range("B43").Value = 0.6800000000000002
$ws.Range("C43").Value = 46.37333333333334
$ws.Range("A44").Value = 19.57306799336643
$ws.Range("B44").Value = 0.8533333333333331
$ws.Range("C44").Value = 54.42000000000001
$ws.Range("A45").Value = 23.85990049751236
$ws.Range("B45").Value = 1.046666666666666
$ws.Range("C45").Value = 63.04666666666665
$ws.Range("A46").Value = 23.98769485903803
$ws.Range("B46").Value = 0.7266666666666666
$ws.Range("C46").Value = 63.86666666666668
$ws.Range("A47").Value = 25.41980099502473
$ws.Range("B47").Value = 0.8933333333333336
$ws.Range("C47").Value = 69.24666666666666
$ws.Range("A48").Value = 18.05489220563847
$ws.Range("B48").Value = 0.6866666666666668
$ws.Range("C48").Value = 47.72666666666665
$ws.Range("A49").Value = 7.512305140961849
$ws.Range("B49").Value = 0.2066666666666667
$ws.Range("C49").Value = 21.65333333333333
$ws.Range("A50").Value = 10.68968490878938
$ws.Range("B50").Value = 0.1200000000000001
$ws.Range("C50").Value = 31.71333333333335
$ws.Range("A51").Value = 4.669718076285232
$ws.Range("B51").Value = 0.03333333333333333
$ws.Range("C51").Value = 14.79333333333333
$ws.Range("A52").Value = 21.93313432835812
$ws.Range("B52").Value = 0.9399999999999997
$ws.Range("C52").Value = 60.50000000000004
$ws.Range("A53").Value = 17.97568822553897
$ws.Range("B53").Value = 0.8933333333333333
$ws.Range("C53").Value = 49.29333333333333
$ws.Range("A54").Value = 7.041592039800984
$ws.Range("B54").Value = 0.2466666666666668
$ws.Range("C54").Value = 21.65999999999999
$ws.Range("A55").Value = 9.402918739635149
$ws.Range("B55").Value = 0.2866666666666667
$ws.Range("C55").Value = 28.32
$ws.Range("A56").Value = 24.50275290215577
$ws.Range("B56").Value = 0.8933333333333334
$ws.Range("C56").Value = 65.14666666666666
$ws.Range("A57").Value = 24.39651741293524
$ws.Range("B57").Value = 0.7666666666666667
$ws.Range("C57").Value = 66.16000000000001
$ws.Range("A58").Value = 11.17595356550579
$ws.Range("B58").Value = 0.46
$ws.Range("C58").Value = 31.7
$ws.Range("A59").Value = 22.88391376451066
$ws.Range("B59").Value = 0.4599999999999998
$ws.Range("C59").Value = 62.66
$ws.Range("A60").Value = 7.980431177446095
$ws.Range("B60").Value = 0.1
$ws.Range("C60").Value = 25.53333333333333
$ws.Range("A61").Value = 23.95794361525699
$ws.Range("B61").Value = 0.9666666666666668
$ws.Range("C61").Value = 66.77333333333335
$ws.Range("A62").Value = 22.14971807628514
$ws.Range("B62").Value = 0.8466666666666665
$ws.Range("C62").Value = 60.05333333333331
$ws.Range("A63").Value = 24.49946932006623
$ws.Range("B63").Value = 0.8200000000000004
$ws.Range("C63").Value = 67.67333333333332
$ws.Range("A64").Value = 7.890646766169143
$ws.Range("B64").Value = 0.1533333333333334
$ws.Range("C64").Value = 23.78666666666667
$ws.Range("A65").Value = 21.90497512437803
$ws.Range("B65").Value = 0.5466666666666667
$ws.Range("C65").Value = 61.31333333333336
$ws.Range("A66").Value = 21.40354892205632
$ws.Range("B66").Value = 0.5666666666666667
$ws.Range("C66").Value = 57.08000000000002
$ws.Range("A67").Value = 19.80447761194025
$ws.Range("B67").Value = 0.7466666666666669
$ws.Range("C67").Value = 55.13333333333333
$ws.Range("A68").Value = 9.054925373134314
$ws.Range("B68").Value = 0.2866666666666668
$ws.Range("C68").Value = 26.00666666666666
$ws.Range("A69").Value = 24.3590381426201
$ws.Range("B69").Value = 0.8733333333333333
$ws.Range("C69").Value = 66.97999999999998
$ws.Range("A70").Value = 8.518507462686557
$ws.Range("B70").Value = 0.1600000000000001
$ws.Range("C70").Value = 25.66
$ws.Range("A71").Value = 12.96606965174128
$ws.Range("B71").Value = 0.3533333333333332
$ws.Range("C71").Value = 36.65333333333333
$ws.Range("A72").Value = 15.36159203980099
$ws.Range("B72").Value = 0.5933333333333333
$ws.Range("C72").Value = 41.77333333333334
